$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in sheet order, EXCLUDING column F (Temp_Diff), which is a
# formula column handled separately below.
# Order: A, B, C, D, E, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$newRows = @(
  @(45835, "Flowering",     "Large",  72, 86, 0.11, 0.2,  "No", 2, "Bright",  8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Nonflowering",  "Medium", 72, 86, 0.11, 0.2,  "No", 3, "Bright",  8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Nonflowering",  "Small",  72, 86, 0.11, 0.15, "No", 3, "Neutral", 8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Nonflowering",  "Medium", 72, 86, 0.11, 0.2,  "No", 3, "Neutral", 8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Nonflowering",  "Medium", 72, 86, 0.11, 0.1,  "No", 3, "Bright",  8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Nonflowering",  "Large",  72, 86, 0.11, 0.4,  "No", 4, "Dark",    8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0),
  @(45835, "Tree",          "Medium", 72, 86, 0.11, 0.95, "No", 1, "Dark",    8, 0.65, 73, 30.04, 13, 0.49, 8.6999999999999993, 50, 0)
)

$colLetters = @("A","B","C","D","E","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$startRow = 338
$r = $startRow
foreach ($rowVals in $newRows) {
  for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $addr = "$($colLetters[$i])$r"
    $ws.Range($addr).Value = $rowVals[$i]
  }
  # Column F: ABS(D-E) formula, continuing the existing shared-formula
  # pattern already used for rows above (F283:F337 -> si="6").
  $ws.Cells.Item($r, 6).Formula = "=ABS(D$r-E$r)"
  $r++
}
$endRow = $r - 1

# Match the original date-number formatting (style) used by column A.
$ws.Range("A337").Copy() | Out-Null
$ws.Range("A$($startRow):A$($endRow)").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reproduce the final selection/view state recorded in the saved workbook.
$ws.Range("R338:R344").Select() | Out-Null

Write-Host "done"
